$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Update existing hours: C64 6 -> 7, C65 3 -> 4
$ws.Cells.Item(64, 3).Value = 7
$ws.Cells.Item(65, 3).Value = 4

# Add a new entry on row 69 (date already present in A69)
$ws.Cells.Item(69, 2).Value = "Implémentation"
$ws.Cells.Item(69, 3).Value = 5
$ws.Cells.Item(69, 4).Value = "Ajout du système central d'accès audevice (accessors) et gestion des accès selon la configuration. KeyboardAccessor"

# Recalculate so the totals row (SUM formula in C77) reflects the new values
$excel.CalculateFullRebuild()

# Update the view so the visible/selected area matches after the edit
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("D72").Select()
